# The Invoice sheet's top section (old rows 1-9: the "INVOICE #100" title,
# company/address/phone placeholders, and the "BILL TO / FOR" block) is
# removed entirely. Deleting those rows shifts the invoice-details table and
# totals block (old rows 10-19) up to become the new rows 1-10, and shrinks
# the InvoiceDetails table/autofilter range from B10:C11 to B1:C2
# automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

$ws.Range("A1:A9").EntireRow.Delete()

# B2 ("COMPANY NAME") - which the Company_Name defined name pointed at - no
# longer exists, so the name now resolves to a broken reference.
$names = $wb.Names
foreach ($n in $names) {
    if ($n.Name -eq "Company_Name") {
        $n.RefersTo = "=Invoice!#REF!"
    }
}

# Leave the sheet's selection on C18, matching the saved view state.
$ws.Application.Goto($ws.Range("C18"))
